$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 414
$ws.Range("I9").Value = 339.6
$ws.Range("K9").Value = 339.6
$ws.Range("M9").Value = -170.6
$ws.Range("H17").Value = 801.78723
$ws.Range("J17").Value = 801.78723
$ws.Range("L17").Value = 2405.36169
$ws.Range("N17").Value = -2741.36169
$ws.Range("H51").Value = 6638.5557
$ws.Range("J51").Value = 6718.4375
$ws.Range("L51").Value = 6718.4375
$ws.Range("N51").Value = -7686.4375
$ws.Range("H76").Value = 12748.5
$ws.Range("I76").Value = 11998
$ws.Range("K76").Value = 11998
$ws.Range("M76").Value = -11683
$ws.Range("H79").Value = 12748.5
$ws.Range("I79").Value = 11998
$ws.Range("K79").Value = 11998
$ws.Range("M79").Value = -10906
$ws.Range("H87").Value = 58145
$ws.Range("J87").Value = 58145
$ws.Range("L87").Value = 58145
$ws.Range("N87").Value = -60641
$ws.Range("H90").Value = 58145
$ws.Range("J90").Value = 58145
$ws.Range("L90").Value = 174435
$ws.Range("N90").Value = -186915
$ws.Range("H92").Value = 188.66667
$ws.Range("I92").Value = 188.66667
$ws.Range("K92").Value = 188.66667
$ws.Range("M92").Value = 1059.33333
$ws.Range("H98").Value = 2116
$ws.Range("I98").Value = 2223.1538
$ws.Range("K98").Value = 2223.1538
$ws.Range("M98").Value = -725.1538
$ws.Range("H103").Value = 226.85715
$ws.Range("I103").Value = 194
$ws.Range("K103").Value = 582
$ws.Range("M103").Value = 4
$ws.Range("H106").Value = 1499
$ws.Range("I106").Value = 1499
$ws.Range("J106").Value = 1499
$ws.Range("K106").Value = 1499
$ws.Range("L106").Value = 1499
$ws.Range("M106").Value = -868
$ws.Range("N106").Value = -2761
$ws.Range("H122").Value = 2116
$ws.Range("I122").Value = 2223.1538
$ws.Range("K122").Value = 6669.4614
$ws.Range("M122").Value = -4219.4614
$ws.Range("H129").Value = 1703.8334
$ws.Range("I129").Value = 1464
$ws.Range("J129").Value = 1943.6666
$ws.Range("K129").Value = 4392
$ws.Range("L129").Value = 5830.9998
$ws.Range("M129").Value = 608
$ws.Range("N129").Value = -15830.9998
$ws.Range("H131").Value = 148239.14
$ws.Range("I131").Value = 148239.14
$ws.Range("K131").Value = 444717.42
$ws.Range("M131").Value = -439677.42
$ws.Range("H137").Value = 1555
$ws.Range("I137").Value = 1101.25
$ws.Range("J137").Value = 2462.5
$ws.Range("K137").Value = 3303.75
$ws.Range("L137").Value = 7387.5
$ws.Range("M137").Value = -753.75
$ws.Range("N137").Value = -12487.5
$ws.Range("H138").Value = 2764.7368
$ws.Range("I138").Value = 3283.889
$ws.Range("J138").Value = 2297.5
$ws.Range("K138").Value = 9851.667000000001
$ws.Range("L138").Value = 6892.5
$ws.Range("M138").Value = -4711.667000000001
$ws.Range("N138").Value = -17172.5
$ws.Range("H141").Value = 2642.375
$ws.Range("I141").Value = 2611.6428
$ws.Range("K141").Value = 7834.928400000001
$ws.Range("M141").Value = -2654.928400000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3974.95
$ws.Range("I45").Value = 1954.5454
$ws.Range("K45").Value = 1954.5454
$ws.Range("M45").Value = -1577.5454
$ws.Range("H61").Value = 5198.9556
$ws.Range("I61").Value = 4134.737
$ws.Range("K61").Value = 4134.737
$ws.Range("M61").Value = -3922.737
$ws.Range("H62").Value = 35000
$ws.Range("J62").Value = 35000
$ws.Range("L62").Value = 35000
$ws.Range("N62").Value = -36248
$ws.Range("H65").Value = 35000
$ws.Range("J65").Value = 35000
$ws.Range("L65").Value = 105000
$ws.Range("N65").Value = -111240
$ws.Range("H110").Value = 6037.524
$ws.Range("J110").Value = 7066.9
$ws.Range("L110").Value = 7066.9
$ws.Range("N110").Value = -11156.9
$ws.Range("H132").Value = 3610.875
$ws.Range("I132").Value = 2694.9143
$ws.Range("K132").Value = 8084.742899999999
$ws.Range("M132").Value = -5554.742899999999
$ws.Range("H136").Value = 5198.9556
$ws.Range("I136").Value = 4134.737
$ws.Range("K136").Value = 12404.211
$ws.Range("M136").Value = -9854.210999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H100").Value = 25280
$ws.Range("J100").Value = 25280
$ws.Range("L100").Value = 25280
$ws.Range("N100").Value = -27444
$ws.Range("H105").Value = 2458.7778
$ws.Range("I105").Value = 2141.125
$ws.Range("K105").Value = 2141.125
$ws.Range("M105").Value = -394.125
$ws.Range("H107").Value = 6628.8
$ws.Range("I107").Value = 6534
$ws.Range("K107").Value = 6534
$ws.Range("M107").Value = -4614

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10069.4
$ws.Range("I58").Value = 5139
$ws.Range("J58").Value = 14999.8
$ws.Range("K58").Value = 5139
$ws.Range("L58").Value = 14999.8
$ws.Range("M58").Value = -4936
$ws.Range("N58").Value = -15405.8
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H107").Value = 530.9167
$ws.Range("J107").Value = 565.4167
$ws.Range("L107").Value = 565.4167
$ws.Range("N107").Value = -4405.4167
$ws.Range("H122").Value = 4245.933
$ws.Range("I122").Value = 4069
$ws.Range("J122").Value = 4599.8
$ws.Range("K122").Value = 12207
$ws.Range("L122").Value = 13799.4
$ws.Range("M122").Value = -9757
$ws.Range("N122").Value = -18699.4
$ws.Range("H134").Value = 3981.0881
$ws.Range("I134").Value = 2918.9167
$ws.Range("J134").Value = 6530.3
$ws.Range("K134").Value = 8756.750100000001
$ws.Range("L134").Value = 19590.9
$ws.Range("M134").Value = -6221.750100000001
$ws.Range("N134").Value = -24660.9
$ws.Range("H136").Value = 10069.4
$ws.Range("I136").Value = 5139
$ws.Range("J136").Value = 14999.8
$ws.Range("K136").Value = 15417
$ws.Range("L136").Value = 44999.39999999999
$ws.Range("M136").Value = -12867
$ws.Range("N136").Value = -50099.39999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 574.9167
$ws.Range("I8").Value = 574.9167
$ws.Range("K8").Value = 1724.7501
$ws.Range("M8").Value = -1585.7501
$ws.Range("H12").Value = 985.1111
$ws.Range("J12").Value = 1052
$ws.Range("L12").Value = 3156
$ws.Range("N12").Value = -3502
$ws.Range("H50").Value = 2249.5
$ws.Range("I50").Value = 2499
$ws.Range("J50").Value = 2000
$ws.Range("K50").Value = 7497
$ws.Range("L50").Value = 6000
$ws.Range("M50").Value = -7016
$ws.Range("N50").Value = -6962
$ws.Range("H53").Value = 2249.5
$ws.Range("I53").Value = 2499
$ws.Range("J53").Value = 2000
$ws.Range("K53").Value = 7497
$ws.Range("L53").Value = 6000
$ws.Range("M53").Value = -7016
$ws.Range("N53").Value = -6962
$ws.Range("H137").Value = 4350.65
$ws.Range("I137").Value = 1070.8948
$ws.Range("K137").Value = 3212.6844
$ws.Range("M137").Value = 1887.3156
$ws.Range("H140").Value = 1773.52
$ws.Range("I140").Value = 1157.8334
$ws.Range("J140").Value = 2341.8462
$ws.Range("K140").Value = 3473.5002
$ws.Range("L140").Value = 7025.5386
$ws.Range("M140").Value = 1706.4998
$ws.Range("N140").Value = -17385.5386

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9355.046
$ws.Range("I46").Value = 3665.52
$ws.Range("J46").Value = 16841.264
$ws.Range("K46").Value = 3665.52
$ws.Range("L46").Value = 16841.264
$ws.Range("M46").Value = -3477.52
$ws.Range("N46").Value = -17217.264
$ws.Range("H103").Value = 33000
$ws.Range("J103").Value = 33000
$ws.Range("L103").Value = 33000
$ws.Range("N103").Value = -35344

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2283
$ws.Range("H4").Value = 16670000
$ws.Range("I4").Value = 16670000
$ws.Range("K4").Value = 16670000
$ws.Range("M4").Value = -16669887
$ws.Range("H5").Value = 79696
$ws.Range("J5").Value = 79696
$ws.Range("L5").Value = 79696
$ws.Range("N5").Value = -79920
$ws.Range("H6").Value = 888888900
$ws.Range("I6").Value = 888888900
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 888888900
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -888888785
$ws.Range("N6").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 2325.3513
$ws.Range("I132").Value = 2224.2942
$ws.Range("K132").Value = 6672.882599999999
$ws.Range("M132").Value = -4142.882599999999
$ws.Range("H136").Value = 5169.1
$ws.Range("I136").Value = 4375.4116
$ws.Range("K136").Value = 13126.2348
$ws.Range("M136").Value = -10576.2348
